$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column F ("On Hand"), shifting it to column H
$ws.Range("F:G").Insert()

# New header cells
$ws.Range("F10").Value = "Color"
$ws.Range("G10").Value = "Sparkling"

# New data cells
$ws.Range("F11").Value = "Red"
$ws.Range("G11").Value = 0

$ws.Range("F12").Value = "White"
$ws.Range("G12").Value = "Yes"

# Match the saved selection / page setup state
$ws.PageSetup.Orientation = 1
$null = $ws.Range("F13").Select()
